# Adds the newly placed order (Order #14, Mrunal) to the "All Orders" log
# and updates the "Daily Summary" sheet with a new 2026-01-14 summary row.
# Both sheets keep their history: the new record is inserted at row 2 and
# every existing row shifts down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "All Orders": insert new row 2 for order #14
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("All Orders")
$ws1.Rows.Item(2).Insert()

$ws1.Cells.Item(2, 1).Value2 = 14                          # Order ID (number)
$ws1.Cells.Item(2, 2).Value2 = "2026-01-14 17:08"          # Date
$ws1.Cells.Item(2, 3).Value2 = "Mrunal"                    # Customer
$ws1.Cells.Item(2, 4).Value2 = "KLV B 2108"                # Flat No

# Phone looks like a pure number - force text so it is not coerced
$ws1.Range("E2").NumberFormat = "@"
$ws1.Cells.Item(2, 5).Value2 = "9404665203"                # Phone
$ws1.Range("E2").Style = "Normal"

$ws1.Cells.Item(2, 6).Value2 = "Wheat Chapati x40"         # Items
$ws1.Cells.Item(2, 7).Value2 = 600                          # Total (number)
$ws1.Cells.Item(2, 8).Value2 = "NEW"                        # Status
$ws1.Cells.Item(2, 9).Value2 = "PENDING"                    # Payment

# Collection Date looks like a date - force text so it is not coerced
$ws1.Range("J2").NumberFormat = "@"
$ws1.Cells.Item(2, 10).Value2 = "2026-01-15"                # Collection Date
$ws1.Range("J2").Style = "Normal"

$ws1.Cells.Item(2, 11).Value2 = "00:30"                     # Collection Time
$ws1.Cells.Item(2, 12).Value2 = ""                          # Notes
$ws1.Cells.Item(2, 13).Value2 = ""                          # Cancel Reason
$ws1.Cells.Item(2, 14).Value2 = ""                          # Feedback

# ---------------------------------------------------------------------
# Sheet "Daily Summary": insert new row 2 for 2026-01-14
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Rows.Item(2).Insert()

# Date column is text formatted (yyyy-mm-dd) - force text so it is not coerced
$ws2.Range("A2").NumberFormat = "@"
$ws2.Cells.Item(2, 1).Value2 = "2026-01-14"                 # Date
$ws2.Range("A2").Style = "Normal"

$ws2.Cells.Item(2, 2).Value2 = 1                            # Total Orders
$ws2.Cells.Item(2, 3).Value2 = 0                            # Delivered
$ws2.Cells.Item(2, 4).Value2 = 0                            # Cancelled
$ws2.Cells.Item(2, 5).Value2 = 600                          # Revenue
$ws2.Cells.Item(2, 6).Value2 = 0                            # Paid
$ws2.Cells.Item(2, 7).Value2 = 600                          # Pending
